$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "chat_id" column before the current column C (after sprite_name),
# shifting all the following columns (max_hp .. Next Group or END) one to the right.
$ws.Columns("C:C").Insert()

# Match the width of the sprite_name column (B) that chat_id (C) now sits next to.
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

$ws.Range("C1").Value = "chat_id"
$ws.Range("C2").Value = "enemy_general_1"
$ws.Range("C3").Value = "tanooki_1"
$ws.Range("C4").Value = "enemy_general_1"
$ws.Range("C5").Value = "enemy_general_1"
$ws.Range("C6").Value = "enemy_general_1"

$ws.Range("C7").Select() | Out-Null
